$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 58 (shifts existing rows 58-109 down to 60-111)
$ws.Range("A58:A59").EntireRow.Insert()

# New row 58 data
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 45225
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 300000000
$ws.Range("G58").Value = "Espárragos"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 800
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = 1500
$ws.Range("N58").Value = "$/kilo"
$ws.Range("O58").Value = "Región de La Araucanía"
$ws.Range("P58").Value = 1500
$ws.Range("Q58").Value = 1
$ws.Range("R58").Value = "Hortaliza"

# New row 59 data
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 45225
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = 300000000
$ws.Range("G59").Value = "Espárragos"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 1400
$ws.Range("K59").Value = 1500
$ws.Range("L59").Value = 1600
$ws.Range("M59").Value = 1571
$ws.Range("N59").Value = "$/kilo"
$ws.Range("O59").Value = "Región del Maule"
$ws.Range("P59").Value = 1571
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"
